# TestData.xlsx - LoginTestData sheet update
# Added a new row-8 test case: page testcases listing page / project overview page / shopping page
# Row 8 user becomes jakay34@gmail.com, and the password column now holds the literal
# numeric-looking string "12345678" (entered with a leading apostrophe so Excel keeps it
# as quoted text instead of coercing it to a number).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginTestData")

$ws.Range("A8").Value = "jakay34@gmail.com"
$ws.Range("B8").Value = "'12345678"

# Match the new active-cell selection left behind in the saved file.
$ws.Range("D5").Select()
